$d = $word.ActiveDocument

# The headline "Full-Stack Developer" is split into three runs:
#   [Full-][Stack][ Developer]
# with w:proofErr spellStart/spellEnd markers wrapping the "Stack" run
# (Word's spell-checker flags "Full-Stack" as two words). The commit
# changes the title to "JavaScript Developer", i.e. merges the first two
# runs' text into "JavaScript" and drops the now-unneeded proofErr markers,
# while leaving " Developer" as its own (separately formatted-but-identical) run.

# Step 1: remove the word "Stack" - this deletes that whole run, leaving the
# (now adjacent, empty) spellStart/spellEnd proofErr markers behind.
$rng1 = $d.Content
$rng1.Find.Execute("Stack", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Step 2: turn "Full-" into "JavaScript" (entirely inside its own run, so it
# doesn't disturb the proofErr markers left over from step 1).
$rng2 = $d.Content
$rng2.Find.Execute("Full-", $true, $false, $false, $false, $false, $true, 1, $false, "JavaScript", 2)

# Step 3: merge "JavaScript" with " Developer" into a single run - this sweeps
# up (and discards) the orphaned spellStart/spellEnd markers sitting between them.
$rng3 = $d.Content
$rng3.Find.Execute("JavaScript Developer", $true, $false, $false, $false, $false, $true, 1, $false, "JavaScript Developer", 2)

# Step 4: re-split that merged run back into "JavaScript" + " Developer" by
# toggling a character property on just the " Developer" portion (off then
# back on, a no-op in the end) which forces the run boundary back without
# reintroducing any proofErr markers.
$rng4 = $d.Content
$rng4.Find.Execute(" Developer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng4.Font.Bold = $false
$rng4.Font.Bold = $true
